$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 511 (what is currently row 511
# "Chilena(o)" / 44421 moves down to become row 513, and so on for every
# row through the former row 526, which becomes row 528).
$ws.Rows.Item(511).Insert()
$ws.Rows.Item(511).Insert()

# Populate the first new row (511)
$ws.Range("A511").Value = 6
$ws.Range("B511").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C511").Value = "Metropolitana"
$ws.Range("D511").Value = 44509
$ws.Range("E511").Value = 13
$ws.Range("F511").Value = 100112021
$ws.Range("G511").Value = "Ají"
$ws.Range("H511").Value = "Americana (o)"
$ws.Range("I511").Value = "Primera"
$ws.Range("J511").Value = 50
$ws.Range("K511").Value = 30000
$ws.Range("L511").Value = 35000
$ws.Range("M511").Value = 32600
$ws.Range("N511").Value = "`$/saco 25 kilos"
$ws.Range("O511").Value = "Provincia de Limarí"
$ws.Range("P511").Value = 1304
$ws.Range("Q511").Value = 25
$ws.Range("R511").Value = "Hortaliza"

# Populate the second new row (512)
$ws.Range("A512").Value = 6
$ws.Range("B512").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C512").Value = "Metropolitana"
$ws.Range("D512").Value = 44509
$ws.Range("E512").Value = 13
$ws.Range("F512").Value = 100112021
$ws.Range("G512").Value = "Ají"
$ws.Range("H512").Value = "Inferno"
$ws.Range("I512").Value = "Primera"
$ws.Range("J512").Value = 70
$ws.Range("K512").Value = 20000
$ws.Range("L512").Value = 25000
$ws.Range("M512").Value = 22143
$ws.Range("N512").Value = "`$/caja 15 kilos"
$ws.Range("O512").Value = "Provincia de Huasco"
$ws.Range("P512").Value = 1476
$ws.Range("Q512").Value = 15
$ws.Range("R512").Value = "Hortaliza"
